$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.30599947411286
$ws.Range("C2").Value = 14.91255542715784
$ws.Range("D2").Value = 6.003839878556669
$ws.Range("E2").Value = 11.5185608246296
$ws.Range("F2").Value = 47.21339889380268
$ws.Range("I2").Value = 30.92631245887682
$ws.Range("J2").Value = 9.994734767596626
$ws.Range("M2").Value = 19.22609602362527
$ws.Range("B3").Value = 16.95595404734549
$ws.Range("C3").Value = 14.54505901709167
$ws.Range("D3").Value = 6.007956793367978
$ws.Range("E3").Value = 11.54033031234061
$ws.Range("F3").Value = 46.95490457852213
$ws.Range("I3").Value = 30.84837308729013
$ws.Range("J3").Value = 10.01472407255188
$ws.Range("M3").Value = 19.14821252242034
$ws.Range("B4").Value = 16.74316502880957
$ws.Range("C4").Value = 14.31988694452889
$ws.Range("D4").Value = 6.010771525705295
$ws.Range("E4").Value = 11.55504241591286
$ws.Range("F4").Value = 46.80737778198909
$ws.Range("I4").Value = 30.80705852837144
$ws.Range("J4").Value = 10.02803715559503
$ws.Range("M4").Value = 19.10509357052726
$ws.Range("B5").Value = 16.65712767487063
$ws.Range("C5").Value = 14.22840015692857
$ws.Range("D5").Value = 6.01199084953917
$ws.Range("E5").Value = 11.56137640673837
$ws.Range("F5").Value = 46.75010596176755
$ws.Range("I5").Value = 30.79186855151262
$ws.Range("J5").Value = 10.03372386056413
$ws.Range("M5").Value = 19.08871772725796
$ws.Range("B6").Value = 16.64288610416648
$ws.Range("C6").Value = 14.21322989339755
$ws.Range("D6").Value = 6.012197688344611
$ws.Range("E6").Value = 11.5624486269367
$ws.Range("F6").Value = 46.74076880942176
$ws.Range("I6").Value = 30.78944569143547
$ws.Range("J6").Value = 10.03468393167074
$ws.Range("M6").Value = 19.08607108948752
$ws.Range("B7").Value = 16.74200177964925
$ws.Range("C7").Value = 14.31865181120753
$ws.Range("D7").Value = 6.010787676987773
$ws.Range("E7").Value = 11.55512646654179
$ws.Range("F7").Value = 46.8065938291146
$ws.Range("I7").Value = 30.80684700591408
$ws.Range("J7").Value = 10.0281127894854
$ws.Range("M7").Value = 19.10486786354643
$ws.Range("B8").Value = 17.18494178462178
$ws.Range("C8").Value = 14.78583557908057
$ws.Range("D8").Value = 6.005199924541205
$ws.Range("E8").Value = 11.52578788036115
$ws.Range("F8").Value = 47.12196869601311
$ws.Range("I8").Value = 30.89808278264799
$ws.Range("J8").Value = 10.00141138196097
$ws.Range("M8").Value = 19.19827428401733
$ws.Range("B9").Value = 18.06457266199407
$ws.Range("C9").Value = 15.69922619141917
$ws.Range("D9").Value = 5.996513010069912
$ws.Range("E9").Value = 11.47891655128442
$ws.Range("F9").Value = 47.82755166858742
$ws.Range("I9").Value = 31.12876371614923
$ws.Range("J9").Value = 9.957293919893145
$ws.Range("M9").Value = 19.41810355754414
$ws.Range("B10").Value = 18.70981898364645
$ws.Range("C10").Value = 16.36037258554958
$ws.Range("D10").Value = 5.991507554871696
$ws.Range("E10").Value = 11.45096006214995
$ws.Range("F10").Value = 48.39656286102616
$ws.Range("I10").Value = 31.32953448924824
$ws.Range("J10").Value = 9.929900547935738
$ws.Range("M10").Value = 19.60099856501484
$ws.Range("B11").Value = 19.001593108003
$ws.Range("C11").Value = 16.65740348452559
$ws.Range("D11").Value = 5.989527976062408
$ws.Range("E11").Value = 11.43964472229577
$ws.Range("F11").Value = 48.66579617806121
$ws.Range("I11").Value = 31.42757284444667
$ws.Range("J11").Value = 9.918527386380942
$ws.Range("M11").Value = 19.68861316026597
$ws.Range("B12").Value = 19.11170809485695
$ws.Range("C12").Value = 16.7692249816392
$ws.Range("D12").Value = 5.988821019680985
$ws.Range("E12").Value = 11.43556117523558
$ws.Range("F12").Value = 48.76918395575381
$ws.Range("I12").Value = 31.46565123484238
$ws.Range("J12").Value = 9.914377070873231
$ws.Range("M12").Value = 19.72240384149362
$ws.Range("B13").Value = 19.08801133884704
$ws.Range("C13").Value = 16.74517330884041
$ws.Range("D13").Value = 5.98897137932624
$ws.Range("E13").Value = 11.4364316911143
$ws.Range("F13").Value = 48.74685469515271
$ws.Range("I13").Value = 31.45740816207829
$ws.Range("J13").Value = 9.915263957717478
$ws.Range("M13").Value = 19.71509951321259
$ws.Range("B14").Value = 19.01066049279162
$ws.Range("C14").Value = 16.6666169504689
$ws.Range("D14").Value = 5.989468959786437
$ws.Range("E14").Value = 11.43930473333627
$ws.Range("F14").Value = 48.67427354125821
$ws.Range("I14").Value = 31.4306865379181
$ws.Range("J14").Value = 9.918182802027713
$ws.Range("M14").Value = 19.69138098503826
$ws.Range("B15").Value = 18.96322872728169
$ws.Range("C15").Value = 16.6184098151068
$ws.Range("D15").Value = 5.98977929577261
$ws.Range("E15").Value = 11.44109076400672
$ws.Range("F15").Value = 48.63000056772486
$ws.Range("I15").Value = 31.41444257713609
$ws.Range("J15").Value = 9.919991052480723
$ws.Range("M15").Value = 19.67693185572226
$ws.Range("B16").Value = 18.69070505486488
$ws.Range("C16").Value = 16.34087529875226
$ws.Range("D16").Value = 5.991642901098184
$ws.Range("E16").Value = 11.45172773597139
$ws.Range("F16").Value = 48.37917190235579
$ws.Range("I16").Value = 31.32326150249403
$ws.Range("J16").Value = 9.930665705895949
$ws.Range("M16").Value = 19.59535973394294
$ws.Range("B17").Value = 18.52298325354569
$ws.Range("C17").Value = 16.1695726288313
$ws.Range("D17").Value = 5.992862266029611
$ws.Range("E17").Value = 11.45861209519101
$ws.Range("F17").Value = 48.22791631494084
$ws.Range("I17").Value = 31.26903528289897
$ws.Range("J17").Value = 9.937492956927546
$ws.Range("M17").Value = 19.54643336957701
$ws.Range("B18").Value = 18.42635609405398
$ws.Range("C18").Value = 16.07069966837092
$ws.Range("D18").Value = 5.993591613551974
$ws.Range("E18").Value = 11.46270379792853
$ws.Range("F18").Value = 48.1418995293949
$ws.Range("I18").Value = 31.238477827571
$ws.Range("J18").Value = 9.941522241404613
$ws.Range("M18").Value = 19.51870913882519
$ws.Range("B19").Value = 18.39361647727075
$ws.Range("C19").Value = 16.03716772240392
$ws.Range("D19").Value = 5.993843370533696
$ws.Range("E19").Value = 11.46411185881628
$ws.Range("F19").Value = 48.11294602146454
$ws.Range("I19").Value = 31.22824045998416
$ws.Range("J19").Value = 9.942904081347351
$ws.Range("M19").Value = 19.50939444626652
$ws.Range("B20").Value = 18.54085478268416
$ws.Range("C20").Value = 16.18784462344619
$ws.Range("D20").Value = 5.99272956531018
$ws.Range("E20").Value = 11.45786558397375
$ws.Range("F20").Value = 48.2439165717137
$ws.Range("I20").Value = 31.27474241846828
$ws.Range("J20").Value = 9.936755583489674
$ws.Range("M20").Value = 19.55159867014412
$ws.Range("B21").Value = 19.03339138430225
$ws.Range("C21").Value = 16.68970962904919
$ws.Range("D21").Value = 5.989321651157202
$ws.Range("E21").Value = 11.43845538975026
$ws.Range("F21").Value = 48.69555392428816
$ws.Range("I21").Value = 31.43850954757173
$ws.Range("J21").Value = 9.917321220822467
$ws.Range("M21").Value = 19.69833123662763
$ws.Range("B22").Value = 19.35306377576029
$ws.Range("C22").Value = 17.01382421664821
$ws.Range("D22").Value = 5.987343035931833
$ws.Range("E22").Value = 11.42694304197184
$ws.Range("F22").Value = 48.99905977435387
$ws.Range("I22").Value = 31.55109171891145
$ws.Range("J22").Value = 9.905531661524719
$ws.Range("M22").Value = 19.79779122209799
$ws.Range("B23").Value = 19.18269091916525
$ws.Range("C23").Value = 16.84123125621267
$ws.Range("D23").Value = 5.988376340561181
$ws.Range("E23").Value = 11.432980142479
$ws.Range("F23").Value = 48.83633053885043
$ws.Range("I23").Value = 31.49050063726109
$ws.Range("J23").Value = 9.911740540189882
$ws.Range("M23").Value = 19.74438915839025
$ws.Range("B24").Value = 18.53277568318602
$ws.Range("C24").Value = 16.17958505629766
$ws.Range("D24").Value = 5.992789471092077
$ws.Range("E24").Value = 11.45820266493632
$ws.Range("F24").Value = 48.23667991688513
$ws.Range("I24").Value = 31.27216029697548
$ws.Range("J24").Value = 9.937088625437219
$ws.Range("M24").Value = 19.54926217715524
$ws.Range("B25").Value = 17.82624497785585
$ws.Range("C25").Value = 15.45330496409088
$ws.Range("D25").Value = 5.998620752132283
$ws.Range("E25").Value = 11.4904571035523
$ws.Range("F25").Value = 47.62758600426041
$ws.Range("I25").Value = 31.06083744745485
$ws.Range("J25").Value = 9.968346919989388
$ws.Range("M25").Value = 19.35480307913832
